# Applies the coin-price-ticker refresh captured in the commit diff:
# updated Price (D) / Volume(1h) (E) figures, and a Filecoin/Stacks row swap
# (rows 31-32: rank order changed, so the two coins traded places).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that Excel would otherwise auto-parse as a Number
# (e.g. "605.22") while preserving the original inline-string/text cell type,
# then drop the transient text NumberFormat so the cell keeps its default style.
function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

$updates = @(
    @{ Cell = "D2"; Value = "70.249.81"; Force = $false },
    @{ Cell = "E2"; Value = "  -0.87%  "; Force = $false },
    @{ Cell = "D3"; Value = "3.510.29"; Force = $false },
    @{ Cell = "E3"; Value = "  -1.15%  "; Force = $false },
    @{ Cell = "E4"; Value = "  +0.01%  "; Force = $false },
    @{ Cell = "D5"; Value = "605.22"; Force = $true },
    @{ Cell = "E5"; Value = "  -0.40%  "; Force = $false },
    @{ Cell = "D6"; Value = "172.84"; Force = $true },
    @{ Cell = "E6"; Value = "  -1.17%  "; Force = $false },
    @{ Cell = "D7"; Value = "0.609"; Force = $true },
    @{ Cell = "E7"; Value = "  -1.60%  "; Force = $false },
    @{ Cell = "D8"; Value = "3.504.57"; Force = $false },
    @{ Cell = "E8"; Value = "  -1.25%  "; Force = $false },
    @{ Cell = "E9"; Value = "  +0.02%  "; Force = $false },
    @{ Cell = "D10"; Value = "0.196"; Force = $true },
    @{ Cell = "E10"; Value = "  -2.45%  "; Force = $false },
    @{ Cell = "D11"; Value = "7.23"; Force = $true },
    @{ Cell = "E11"; Value = "  +6.87%  "; Force = $false },
    @{ Cell = "D12"; Value = "0.587"; Force = $true },
    @{ Cell = "E12"; Value = "  -0.40%  "; Force = $false },
    @{ Cell = "D13"; Value = "46.25"; Force = $true },
    @{ Cell = "E13"; Value = "  -3.45%  "; Force = $false },
    @{ Cell = "E14"; Value = "  -1.83%  "; Force = $false },
    @{ Cell = "D15"; Value = "4.080.66"; Force = $false },
    @{ Cell = "E15"; Value = "  -0.96%  "; Force = $false },
    @{ Cell = "E16"; Value = "  -1.04%  "; Force = $false },
    @{ Cell = "D17"; Value = "612.89"; Force = $true },
    @{ Cell = "E17"; Value = "  -2.89%  "; Force = $false },
    @{ Cell = "D18"; Value = "3.504.58"; Force = $false },
    @{ Cell = "E18"; Value = "  -1.27%  "; Force = $false },
    @{ Cell = "D19"; Value = "70.245.58"; Force = $false },
    @{ Cell = "E19"; Value = "  -0.86%  "; Force = $false },
    @{ Cell = "E20"; Value = "  +0.70%  "; Force = $false },
    @{ Cell = "D21"; Value = "17.52"; Force = $true },
    @{ Cell = "E21"; Value = "  +0.12%  "; Force = $false },
    @{ Cell = "D22"; Value = "0.879"; Force = $true },
    @{ Cell = "E22"; Value = "  -1.51%  "; Force = $false },
    @{ Cell = "E23"; Value = "  -9.62%  "; Force = $false },
    @{ Cell = "D24"; Value = "98.69"; Force = $true },
    @{ Cell = "E24"; Value = "  +1.44%  "; Force = $false },
    @{ Cell = "E25"; Value = "  -2.19%  "; Force = $false },
    @{ Cell = "E26"; Value = "  -3.91%  "; Force = $false },
    @{ Cell = "E27"; Value = "  -0.03%  "; Force = $false },
    @{ Cell = "D28"; Value = "2.56"; Force = $true },
    @{ Cell = "E28"; Value = "  -2.75%  "; Force = $false },
    @{ Cell = "D29"; Value = "33.96"; Force = $true },
    @{ Cell = "E29"; Value = "  +1.37%  "; Force = $false },
    @{ Cell = "E30"; Value = "  -3.30%  "; Force = $false },
    @{ Cell = "B31"; Value = "Stacks"; Force = $false },
    @{ Cell = "C31"; Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; Force = $false },
    @{ Cell = "D31"; Value = "2.98"; Force = $true },
    @{ Cell = "E31"; Value = "  -4.91%  "; Force = $false },
    @{ Cell = "B32"; Value = "Filecoin"; Force = $false },
    @{ Cell = "C32"; Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; Force = $false },
    @{ Cell = "D32"; Value = "8.06"; Force = $true },
    @{ Cell = "E32"; Value = "  -5.19%  "; Force = $false },
    @{ Cell = "E33"; Value = "  -5.31%  "; Force = $false },
    @{ Cell = "E34"; Value = "  -3.82%  "; Force = $false },
    @{ Cell = "D35"; Value = "631.29"; Force = $true },
    @{ Cell = "E35"; Value = "  +10.59%  "; Force = $false },
    @{ Cell = "D36"; Value = "0.0997"; Force = $true },
    @{ Cell = "E36"; Value = "  -2.80%  "; Force = $false },
    @{ Cell = "E37"; Value = "  -0.71%  "; Force = $false },
    @{ Cell = "D38"; Value = "0.0485"; Force = $true },
    @{ Cell = "E38"; Value = "  +6.41%  "; Force = $false },
    @{ Cell = "E39"; Value = "  -4.98%  "; Force = $false },
    @{ Cell = "D40"; Value = "56.85"; Force = $true },
    @{ Cell = "E40"; Value = "  -1.38%  "; Force = $false },
    @{ Cell = "E41"; Value = "  +0.02%  "; Force = $false },
    @{ Cell = "D42"; Value = "0.145"; Force = $true },
    @{ Cell = "E42"; Value = "  +0.95%  "; Force = $false },
    @{ Cell = "D43"; Value = "3.373.56"; Force = $false },
    @{ Cell = "E43"; Value = "  +0.69%  "; Force = $false },
    @{ Cell = "D44"; Value = "0.0₃0735"; Force = $false },
    @{ Cell = "E44"; Value = "  +1.76%  "; Force = $false },
    @{ Cell = "D45"; Value = "0.311"; Force = $true },
    @{ Cell = "E45"; Value = "  -6.22%  "; Force = $false },
    @{ Cell = "E46"; Value = "  -4.87%  "; Force = $false },
    @{ Cell = "D47"; Value = "31.95"; Force = $true },
    @{ Cell = "E47"; Value = "  -4.07%  "; Force = $false },
    @{ Cell = "E48"; Value = "  -4.76%  "; Force = $false },
    @{ Cell = "E49"; Value = "  +0.12%  "; Force = $false },
    @{ Cell = "D50"; Value = "133.06"; Force = $true }
)

foreach ($u in $updates) {
    if ($u.Force) {
        Set-TextValue $u.Cell $u.Value
    } else {
        $ws.Range($u.Cell).Value = $u.Value
    }
}
